$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds numeric-looking text (e.g. "0.9997", "9.300",
# "29.967.85") that must stay exactly as text and not be auto-converted to
# a number by Excel. Format the range as Text first, write all the updated
# values, then restore the default cell style so no stray formatting is
# left behind.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "29.967.85"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "1.883.21"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "0.7431"
$ws.Range("E5").Value = "  -3.29%  "
$ws.Range("D6").Value = "243.06"
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("D7").Value = "0.9998"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "0.3161"
$ws.Range("E8").Value = "  +0.78%  "
$ws.Range("D9").Value = "0.07220"
$ws.Range("E9").Value = "  +1.23%  "
$ws.Range("D10").Value = "24.93"
$ws.Range("E10").Value = "  -2.84%  "
$ws.Range("D11").Value = "0.08340"
$ws.Range("E11").Value = "  -2.23%  "
$ws.Range("D12").Value = "1.924.81"
$ws.Range("E12").Value = "  +1.37%  "
$ws.Range("D13").Value = "0.7561"
$ws.Range("E13").Value = "  -1.06%  "
$ws.Range("D14").Value = "5.415"
$ws.Range("E14").Value = "  +0.82%  "
$ws.Range("D15").Value = "92.60"
$ws.Range("E15").Value = "  -1.20%  "
$ws.Range("D16").Value = "6.173"
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("D17").Value = "29.935.30"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("D18").Value = "249.83"
$ws.Range("E18").Value = "  +2.23%  "
$ws.Range("D19").Value = "13.61"
$ws.Range("E19").Value = "  -1.14%  "
$ws.Range("D20").Value = "0.000007853"
$ws.Range("E20").Value = "  +0.44%  "
$ws.Range("D21").Value = "2.219.27"
$ws.Range("E21").Value = "  +3.42%  "
$ws.Range("D22").Value = "0.9986"
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").Value = "8.011"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").Value = "0.9999"
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").Value = "0.1561"
$ws.Range("E25").Value = "  -4.16%  "
$ws.Range("D26").Value = "9.300"
$ws.Range("E26").Value = "  -0.90%  "
$ws.Range("D27").Value = "165.53"
$ws.Range("E27").Value = "  +1.43%  "
$ws.Range("D28").Value = "18.73"
$ws.Range("E28").Value = "  -0.32%  "
$ws.Range("D29").Value = "2.042"
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").Value = "1.489"
$ws.Range("E30").Value = "  -1.93%  "
$ws.Range("D31").Value = "4.616"
$ws.Range("E31").Value = "  +2.38%  "
$ws.Range("E32").Value = "  +0.22%  "
$ws.Range("D33").Value = "4.236"
$ws.Range("E33").Value = "  +2.84%  "
$ws.Range("D34").Value = "0.05374"
$ws.Range("E34").Value = "  -1.45%  "
$ws.Range("E35").Value = "  +0.89%  "
$ws.Range("D36").Value = "0.7577"
$ws.Range("E36").Value = "  +1.58%  "
$ws.Range("D37").Value = "0.9953"
$ws.Range("E37").Value = "  -0.56%  "
$ws.Range("D38").Value = "2.707"
$ws.Range("E38").Value = "  +0.15%  "
$ws.Range("D39").Value = "0.01968"
$ws.Range("E39").Value = "  +0.98%  "
$ws.Range("E40").Value = "  -0.59%  "
$ws.Range("D41").Value = "0.4565"
$ws.Range("E41").Value = "  +2.00%  "
$ws.Range("D42").Value = "1.105.39"
$ws.Range("E42").Value = "  +0.30%  "
$ws.Range("D43").Value = "6.066"
$ws.Range("E43").Value = "  -0.29%  "
$ws.Range("D44").Value = "72.77"
$ws.Range("E44").Value = "  -0.63%  "
$ws.Range("D45").Value = "0.8717"
$ws.Range("E45").Value = "  +1.96%  "
$ws.Range("D46").Value = "104.39"
$ws.Range("E46").Value = "  +1.25%  "
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").Value = "1.869"
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("D49").Value = "7.619"
$ws.Range("E49").Value = "  -0.90%  "
$ws.Range("D50").Value = "2.064.59"
$ws.Range("E50").Value = "  +1.81%  "
$ws.Range("D51").Value = "9.556"
$ws.Range("E51").Value = "  -1.99%  "

$priceRange.Style = "Normal"

Write-Host "Updated crypto price table"
